# HeuristicLab 3.3 - Architecture.pptx
# #2330 - deleted overview slides - updated architecture slides and the
# tutorials that used that chart.
#
# This script reproduces (as closely as the PowerPoint COM-interop surface
# allows) the three content edits on the single remaining "architecture"
# slide plus the date-stamp refresh that ripples through the slide master
# and every slide layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder ("18.02.2014" -> "09.03.2015") on the slide master
#    and on every one of its custom layouts. PowerPoint re-stamps this
#    automatic date field whenever the deck is saved on a different day;
#    here we simply rewrite the cached text wherever we find it.
# ---------------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "18.02.2014") {
                $tr.Text = "09.03.2015"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholders $layouts.Item($L).Shapes
}

# ---------------------------------------------------------------------
# 2) Architecture slide edits
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if (-not $shape.HasTextFrame) {
        continue
    }
    $tr = $shape.TextFrame.TextRange
    $text = $tr.Text

    # "Operators" -> "Operators.*" (single run, formatting untouched)
    if ($text -eq "Operators") {
        $tr.Text = "Operators.*"
    }

    # "Programmable Op." -> "Scripting"
    elseif ($text -eq "Programmable Op.") {
        # Fold the " Op." run's text into "Scripting" (keeps that run's
        # clean formatting, i.e. no spell-check "err" flag), then drop
        # the leading "Programmable" run entirely.
        $secondRun = $tr.Characters(13, 4)
        $secondRun.Text = "Scripting"
        $firstRun = $tr.Characters(1, 12)
        $firstRun.Text = ""
    }

    # "Microsoft .NET 4.0" -> "Microsoft .NET " + "4.5" (two runs)
    elseif ($text -eq "Microsoft .NET 4.0") {
        $verRun = $tr.Characters(16, 3)
        $verRun.Text = "4.5"
    }
}
